# chore: update Sheets via scheduled runner
# Refreshes the market-price-derived columns (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ -> cols H:N)
# for the rows whose underlying market data changed, across the ALC, ARM,
# BSM, CRP, CUL, GSM, LTW and WVR sheets. A few rows also lose their
# NQ-or-HQ profit cell entirely (no HQ recipe priced that refresh), so
# those are cleared instead of being set.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 686.84
$ws.Range("I53").Value = 77.416664
$ws.Range("J53").Value = 1249.3846
$ws.Range("K53").Value = 77.416664
$ws.Range("L53").Value = 1249.3846
$ws.Range("M53").Value = 559.583336
$ws.Range("N53").Value = -2523.3846

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 814.3333
$ws.Range("I82").Value = 814.3333
$ws.Range("K82").Value = 2442.9999
$ws.Range("M82").Value = -2036.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H85").Value = 814.3333
$ws.Range("I85").Value = 814.3333
$ws.Range("K85").Value = 2442.9999
$ws.Range("M85").Value = -1038.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1707.1143
$ws.Range("I132").Value = 1088.421
$ws.Range("K132").Value = 3265.263
$ws.Range("M132").Value = -735.2629999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 4552697
$ws.Range("I135").Value = 6251340.5
$ws.Range("J135").Value = 22981
$ws.Range("K135").Value = 56262064.5
$ws.Range("L135").Value = 206829
$ws.Range("M135").Value = -56259529.5
$ws.Range("N135").Value = -211899

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 7904.7437
$ws.Range("I137").Value = 5623.8
$ws.Range("J137").Value = 8691.275
$ws.Range("K137").Value = 16871.4
$ws.Range("L137").Value = 26073.825
$ws.Range("M137").Value = -14321.4
$ws.Range("N137").Value = -31173.825

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4872.9473
$ws.Range("I138").Value = 4781.625
$ws.Range("J138").Value = 4939.364
$ws.Range("K138").Value = 14344.875
$ws.Range("L138").Value = 14818.092
$ws.Range("M138").Value = -9204.875
$ws.Range("N138").Value = -25098.092

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 99988
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2387.3333
$ws.Range("I45").Value = 2387.3333
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2387.3333
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2010.3333
$ws.Range("N45").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4159.7036
$ws.Range("I61").Value = 3596.2727
$ws.Range("K61").Value = 3596.2727
$ws.Range("M61").Value = -3384.2727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 8260.040000000001
$ws.Range("I132").Value = 5794.5293
$ws.Range("K132").Value = 17383.5879
$ws.Range("M132").Value = -14853.5879

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4159.7036
$ws.Range("I136").Value = 3596.2727
$ws.Range("K136").Value = 10788.8181
$ws.Range("M136").Value = -8238.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 112023.22
$ws.Range("I80").Value = 1018.8
$ws.Range("J80").Value = 250778.75
$ws.Range("K80").Value = 1018.8
$ws.Range("L80").Value = 250778.75
$ws.Range("M80").Value = -20.79999999999995
$ws.Range("N80").Value = -252774.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 112023.22
$ws.Range("I83").Value = 1018.8
$ws.Range("J83").Value = 250778.75
$ws.Range("K83").Value = 5094
$ws.Range("L83").Value = 1253893.75
$ws.Range("M83").Value = -102
$ws.Range("N83").Value = -1263877.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3099.9812
$ws.Range("I31").Value = 1313.3667
$ws.Range("K31").Value = 1313.3667
$ws.Range("M31").Value = -1018.3667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3099.9812
$ws.Range("I34").Value = 1313.3667
$ws.Range("K34").Value = 1313.3667
$ws.Range("M34").Value = -1111.3667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7741.386
$ws.Range("I58").Value = 6185.647
$ws.Range("J58").Value = 8720.925999999999
$ws.Range("K58").Value = 6185.647
$ws.Range("L58").Value = 8720.925999999999
$ws.Range("M58").Value = -5982.647
$ws.Range("N58").Value = -9126.925999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 45755.137
$ws.Range("I132").Value = 3322.8
$ws.Range("K132").Value = 9968.400000000001
$ws.Range("M132").Value = -7438.400000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 8341.833000000001
$ws.Range("I134").Value = 8350.25
$ws.Range("K134").Value = 25050.75
$ws.Range("M134").Value = -22515.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 7741.386
$ws.Range("I136").Value = 6185.647
$ws.Range("J136").Value = 8720.925999999999
$ws.Range("K136").Value = 18556.941
$ws.Range("L136").Value = 26162.778
$ws.Range("M136").Value = -16006.941
$ws.Range("N136").Value = -31262.778

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 2611.8667
$ws.Range("I92").Value = 1729
$ws.Range("J92").Value = 3200.4443
$ws.Range("K92").Value = 5187
$ws.Range("L92").Value = 9601.332900000001
$ws.Range("M92").Value = -3939
$ws.Range("N92").Value = -12097.3329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1742.8334
$ws.Range("I102").Value = 1728.0358
$ws.Range("J102").Value = 1950
$ws.Range("K102").Value = 1728.0358
$ws.Range("L102").Value = 1950
$ws.Range("M102").Value = -106.0358000000001
$ws.Range("N102").Value = -5194

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3518.7334
$ws.Range("I132").Value = 2124.5
$ws.Range("J132").Value = 6307.2
$ws.Range("K132").Value = 6373.5
$ws.Range("L132").Value = 18921.6
$ws.Range("M132").Value = -3843.5
$ws.Range("N132").Value = -23981.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8730.299999999999
$ws.Range("I7").Value = 4145.3335
$ws.Range("J7").Value = 49995
$ws.Range("K7").Value = 4145.3335
$ws.Range("L7").Value = 49995
$ws.Range("M7").Value = -4033.3335
$ws.Range("N7").Value = -50219

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6660.28
$ws.Range("I40").Value = 6660.28
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 6660.28
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -6524.28
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 8730.299999999999
$ws.Range("I126").Value = 4145.3335
$ws.Range("J126").Value = 49995
$ws.Range("K126").Value = 12436.0005
$ws.Range("L126").Value = 149985
$ws.Range("M126").Value = -9966.000499999998
$ws.Range("N126").Value = -154925

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 16835816
$ws.Range("I126").Value = 16835816
$ws.Range("K126").Value = 50507448
$ws.Range("M126").Value = -50504978

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6665.359
$ws.Range("I132").Value = 5293.1816
$ws.Range("J132").Value = 8441.117
$ws.Range("K132").Value = 15879.5448
$ws.Range("L132").Value = 25323.351
$ws.Range("M132").Value = -13349.5448
$ws.Range("N132").Value = -30383.351

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3001.926
$ws.Range("I136").Value = 1952.35
$ws.Range("J136").Value = 6000.7144
$ws.Range("K136").Value = 5857.049999999999
$ws.Range("L136").Value = 18002.1432
$ws.Range("M136").Value = -3307.049999999999
$ws.Range("N136").Value = -23102.1432
